$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - add new columns I and J, matching the style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for new columns I and J
$data = @(
    @(7, 8),
    @(4, 6),
    @(1, 3),
    @(1, 3),
    @(1, 4),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
